# Auto-generated COM-interop edit script
# Manually add ICDC TCL01 Test for Prasanna (StudyFilesTab row + updated Cypher
# queries for CasesTab / SamplesTab / FilesTab).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the CasesTab row (row 2): new query text + new stat query text,
#    plus vertical-center/wrap formatting on the stat-query cell.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nMATCH (samp:sample)-->(c)`nWHERE  samp.specific_sample_pathology in ['Pulmonary Adenocarcinoma']`nWITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,`n       coalesce(demo.sex, '') AS Sex,`n       coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`ncoalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS ``Weight (kg)``,`n       coalesce(diag.best_response, '') AS ``Response to Treatment``,`n       coalesce(co.cohort_description, '') AS ``Cohort```norder by c.case_id asc`nlimit 100"
$ws.Range("C2").Value = "MATCH (p:program)<--(s:study)<--(c)`nMATCH (cf)-->(samp:sample)`nWHERE samp.specific_sample_pathology IN ['Pulmonary Adenocarcinoma']`nMATCH (cf:file)-[*]->(c:case)`nOPTIONAL MATCH (sf:file)-->(s)`nRETURN`n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n   count(distinct samp) AS Samples,`n    count(distinct cf) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"
$ws.Range("C2").VerticalAlignment = -4108
$ws.Range("C2").WrapText = $true

# ---------------------------------------------------------------------------
# 2. Update the SamplesTab row (row 3): new query text (18pt font + wrap),
#    new stat query text (vertical-center/wrap), filename swapped to the
#    workbook's own filename for column D.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) `nWHERE  samp.specific_sample_pathology IN ['Pulmonary Adenocarcinoma']`nWITH DISTINCT samp AS samp, c, demo, diag`nRETURN  coalesce(samp.sample_id, '') AS ``Sample ID``, `n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(demo.breed,'') AS Breed,`n        coalesce(diag.disease_term,'') AS Diagnosis, `n        coalesce(samp.sample_site, '') AS ``Sample Site``,`n        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,`n        coalesce(samp.specific_sample_pathology, '') AS ``Pathology/Morphology``,`n        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,`n        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,`n        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,`n        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,`n        coalesce(samp.sample_preservation, '') AS ``Sample Preservation```norder by samp.sample_id asc`nlimit 200"
$ws.Range("B3").Font.Size = 18
$ws.Range("B3").WrapText = $true
$ws.Range("C3").Value = "MATCH (p:program)<--(s:study)<--(c)`nMATCH (cf)-->(samp:sample)`nWHERE samp.specific_sample_pathology IN ['Pulmonary Adenocarcinoma']`nMATCH (cf:file)-[*]->(c:case)`nOPTIONAL MATCH (sf:file)-->(s)`nRETURN`n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct cf) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"
$ws.Range("C3").VerticalAlignment = -4108
$ws.Range("C3").WrapText = $true
$ws.Range("D3").Value = "TC08_Canine_Filter_SamplePatho-PulmoAdeno.xlsx"

# ---------------------------------------------------------------------------
# 3. Update the FilesTab row (row 4): new query text (18pt font + wrap),
#    new stat query text (vertical-center/wrap).
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nMATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nMATCH (f)-->(samp:sample)`nWHERE samp.specific_sample_pathology IN ['Pulmonary Adenocarcinoma']`n MATCH (f)-[*]->(samp:sample)`nWITH`n        DISTINCT f, parent, c, demo, diag, s, samp,`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH`n        f, parent, c, demo, diag, s, samp,`n        f.file_size /(1024^i) AS value,`n        10^precision AS factor,`n        units[i] as unit`nWITH`n        f, parent, c, demo, diag, s, samp, unit,`n        round(factor * value)/factor AS size`nRETURN`n        coalesce(f.file_name, '') AS ``File Name``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_type, '') AS ``File Type``,`n        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(samp.sample_id, '') AS ``Sample ID``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(demo.breed,'') AS Breed ,`n        coalesce(diag.disease_term,'') AS Diagnosis`n        order by f.file_name asc`n        limit 200"
$ws.Range("B4").Font.Size = 18
$ws.Range("B4").WrapText = $true
$ws.Range("C4").Value = "MATCH (p:program)<--(s:study)<--(c)`nMATCH (cf)-->(samp:sample)`nWHERE samp.specific_sample_pathology IN ['Pulmonary Adenocarcinoma']`nMATCH (cf:file)-[*]->(c:case)`nOPTIONAL MATCH (sf:file)-->(s)`nRETURN`n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct cf) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"
$ws.Range("C4").VerticalAlignment = -4108
$ws.Range("C4").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Insert the new StudyFilesTab row (row 5) with the same 5 columns.
# ---------------------------------------------------------------------------
$ws.Rows("5:5").Insert()
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = "MATCH (f:file)-->(s:study)`nMATCH (s)<--(c:case)<--(diag:diagnosis)`nMATCH (c)<--(demo:demographic)`nMATCH (samp:sample)-->(c)`nWHERE samp.specific_sample_pathology IN ['Pulmonary Adenocarcinoma']`nWITH`n        DISTINCT f, c, demo, diag, s,`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH`n        f, c, demo, diag, s,`n        f.file_size /(1024^i) AS value, 10^precision AS factor,`n        units[i] as unit`n        WITH`n        f,  c, demo, diag, s, unit,`n        round(factor * value)/factor AS size`nRETURN DISTINCT`n  coalesce(f.file_name, '') AS ``File Name``,`n  coalesce(f.file_type, '') AS ``File Type``,`n  coalesce(`"study`", '') AS ``Association``,`n  coalesce(f.file_description, '') AS ``Description``,`n  coalesce(f.file_format, '') AS  Format,`n  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n  coalesce(s.clinical_study_designation,'') AS ``Study Code```n  order by 'File Name' asc`n  limit 100"
$ws.Range("B5").Font.Size = 18
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = "MATCH (p:program)<--(s:study)<--(c)`nMATCH (cf)-->(samp:sample)`nWHERE samp.specific_sample_pathology IN ['Pulmonary Adenocarcinoma']`nMATCH (cf:file)-[*]->(c:case)`nOPTIONAL MATCH (sf:file)-->(s)`nRETURN`n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct cf) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"
$ws.Range("C5").Font.Size = 18
$ws.Range("C5").WrapText = $true
$ws.Range("C5").VerticalAlignment = -4107
$ws.Range("D5").Value = "TC08_Canine_Filter_SamplePatho-PulmoAdeno_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC08_Canine_Filter_SamplePatho-PulmoAdeno_WebData.xlsx"

# ---------------------------------------------------------------------------
# 5. Row heights (rows grew/shrank once the text + font changed).
# ---------------------------------------------------------------------------
$ws.Rows("2:2").RowHeight = 304.5
$ws.Rows("3:3").RowHeight = 409.5
$ws.Rows("4:4").RowHeight = 409.5
$ws.Rows("5:5").RowHeight = 409.5

# ---------------------------------------------------------------------------
# 6. View state: scroll/zoom/selection to match the saved window position.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 69
$ws.Range("G5").Select()

Write-Output "edit complete"
